$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 (H3-0 DEDM): fill in sample size / vert div / hor div columns
$ws.Range("B6").Value = "4.9"
$ws.Range("C6").Value = "full beam"
$ws.Range("D6").Value = "+-1"
$ws.Range("E6").Value = "+-1"

# Row 9 (H3-3 LIRA): add numeric sample size
$ws.Range("C9").Value = 4

# Row 10 (H3-3 SONATA): add numeric sample size
$ws.Range("C10").Value = 1

# Row 11 (H3-4 SANS-3): sample size text fix 1х1 -> 3x3
$ws.Range("C11").Value = "3x3"

# Row 12 (H3-4 SANS-2): sample size text fix 1х1 -> 3x3
$ws.Range("C12").Value = "3x3"

# Row 14 (H3-4 Tensor): sample size text fix 1х1 -> 5x5
$ws.Range("C14").Value = "5x5"

# Row 15 (H3-5 TOF): add sample size
$ws.Range("C15").Value = "1x1"

# Update selection to C15
$ws.Range("C15").Select() | Out-Null
